$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 592; this shifts rows 592..662 down to 593..663
# (matching the rest of the existing "Start Ruby" data set already on the sheet).
$ws.Rows("592:592").Insert()

# Populate the newly inserted row 592 with the new record.
$ws.Range("A592").Value = 4
$ws.Range("B592").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C592").Value = "Los Lagos"
$ws.Range("D592").Value = 45154
$ws.Range("E592").Value = 10
$ws.Range("F592").Value = "Fruta"
$ws.Range("G592").Value = 100102
$ws.Range("H592").Value = "Cítricos"
$ws.Range("I592").Value = 100102006
$ws.Range("J592").Value = "Pomelo"
$ws.Range("K592").Value = "Start Ruby"
$ws.Range("L592").Value = "Primera"
$ws.Range("M592").Value = 120
$ws.Range("N592").Value = 14000
$ws.Range("O592").Value = 15000
$ws.Range("P592").Value = 14500
$ws.Range("Q592").Value = "$/caja 14 kilos empedrada"
$ws.Range("R592").Value = "Región de O'Higgins"
$ws.Range("S592").Value = 1036
$ws.Range("T592").Value = 14
